$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all review rows (rows 2-21), keeping only row 1 (header/count cell in B1)
$ws.Range("A2:B21").EntireRow.Delete()
